$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the row above (row 41) down to row 42 first so the new row
# inherits the exact same cell styles (date format, fills, borders)
# without minting any new style entries.
$ws.Range("A41:C41").Copy()
$ws.Range("A42:C42").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Add a new log entry row (row 42) for the Django work-log entry.
$ws.Range("A42").Value = 43350
$ws.Range("B42").Value = "python-Django"
$ws.Range("C42").Value = "installation and basics,,started a polls app project"

# Update selection to match the new last cell.
$ws.Range("C42").Select()

